# Update the "Förändrad" (Changed) date column (C) from 2023-10-06 (45205)
# to 2023-10-07 (45206) for every data row (rows 2 through 351).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 351
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
